$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $cell = $ws.Range($cellRef)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}

Set-TextValue "D2" "43.946.24"
Set-TextValue "E2" "  +1.84%  "
Set-TextValue "D3" "2.361.80"
Set-TextValue "E3" "  +0.37%  "
Set-TextValue "E4" "  +0.20%  "
Set-TextValue "D5" "0.690"
Set-TextValue "E5" "  +5.96%  "
Set-TextValue "D6" "241.33"
Set-TextValue "E6" "  +3.22%  "
Set-TextValue "D7" "75.91"
Set-TextValue "E7" "  +5.59%  "
Set-TextValue "E8" "  +0.00%  "
Set-TextValue "D9" "0.627"
Set-TextValue "E9" "  +26.20%  "
Set-TextValue "D10" "0.103"
Set-TextValue "E10" "  +5.56%  "
Set-TextValue "D11" "57.29"
Set-TextValue "E11" "  +0.67%  "
Set-TextValue "D12" "32.88"
Set-TextValue "E12" "  +19.92%  "
Set-TextValue "D13" "7.51"
Set-TextValue "E13" "  +19.22%  "
Set-TextValue "E14" "  +1.37%  "
Set-TextValue "D15" "2.714.45"
Set-TextValue "E15" "  +0.44%  "
Set-TextValue "D16" "16.85"
Set-TextValue "E16" "  +4.08%  "
Set-TextValue "E17" "  +6.16%  "
Set-TextValue "D18" "2.362.08"
Set-TextValue "E18" "  +0.24%  "
Set-TextValue "D19" "43.941.88"
Set-TextValue "E19" "  +1.74%  "
Set-TextValue "E20" "  +3.08%  "
Set-TextValue "D21" "6.68"
Set-TextValue "E21" "  +5.15%  "
Set-TextValue "D22" "77.60"
Set-TextValue "E22" "  +4.19%  "
Set-TextValue "D23" "256.83"
Set-TextValue "E23" "  +2.40%  "
Set-TextValue "E24" "  -0.01%  "
Set-TextValue "E25" "  +2.87%  "
Set-TextValue "D26" "11.24"
Set-TextValue "E26" "  +11.92%  "
Set-TextValue "E27" "  -5.40%  "
Set-TextValue "D28" "1.77"
Set-TextValue "E28" "  +15.21%  "
Set-TextValue "D29" "2.30"
Set-TextValue "E29" "  +1.89%  "
Set-TextValue "D30" "23.19"
Set-TextValue "E30" "  +3.34%  "
Set-TextValue "D31" "175.76"
Set-TextValue "E31" "  +1.90%  "
Set-TextValue "E32" "  -2.74%  "
Set-TextValue "E33" "  +5.35%  "
Set-TextValue "D34" "5.31"
Set-TextValue "E34" "  +6.27%  "
Set-TextValue "D35" "0.0754"
Set-TextValue "E35" "  +9.05%  "
Set-TextValue "E36" "  +6.07%  "
Set-TextValue "D37" "3.83"
Set-TextValue "E37" "  +2.43%  "
Set-TextValue "D38" "2.44"
Set-TextValue "E38" "  +0.51%  "
Set-TextValue "D39" "6.51"
Set-TextValue "E39" "  -0.03%  "
Set-TextValue "D40" "0.0277"
Set-TextValue "E40" "  +8.72%  "
Set-TextValue "E41" "  +19.61%  "
Set-TextValue "E42" "  -0.15%  "
Set-TextValue "B43" "BinanceUSD"
Set-TextValue "C43" "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextValue "D43" "1.00"
Set-TextValue "E43" "  +0.07%  "
Set-TextValue "B44" "FraxShare"
Set-TextValue "C44" "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D44" "8.94"
Set-TextValue "E44" "  +0.69%  "
Set-TextValue "E45" "  +5.57%  "
Set-TextValue "E46" "  +5.03%  "
Set-TextValue "E47" "  +13.17%  "
Set-TextValue "D48" "102.10"
Set-TextValue "E48" "  +2.96%  "
Set-TextValue "E49" "  +2.93%  "
Set-TextValue "D50" "4.48"
Set-TextValue "E50" "  -0.38%  "
Set-TextValue "D51" "54.45"
Set-TextValue "E51" "  +7.32%  "
